$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.581.35'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '2.946.39'
$ws.Range('E3').Value = '  -2.29%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '375.39'
$ws.Range('E5').Value = '  +5.52%  '
$ws.Range('D6').Value = '105.59'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('E7').Value = '  -2.76%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').Value = '0.596'
$ws.Range('E9').Value = '  -4.45%  '
$ws.Range('D10').Value = '37.43'
$ws.Range('E10').Value = '  -3.05%  '
$ws.Range('D11').Value = '0.139'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '0.0841'
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('D13').Value = '18.50'
$ws.Range('E13').Value = '  -4.52%  '
$ws.Range('D14').Value = '3.411.16'
$ws.Range('E14').Value = '  -2.46%  '
$ws.Range('D15').Value = '7.48'
$ws.Range('E15').Value = '  -3.45%  '
$ws.Range('D16').Value = '2.947.81'
$ws.Range('E16').Value = '  -2.64%  '
$ws.Range('D17').Value = '0.944'
$ws.Range('E17').Value = '  -8.03%  '
$ws.Range('D18').Value = '51.591.50'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('E19').Value = '  -6.43%  '
$ws.Range('D20').Value = '7.35'
$ws.Range('E20').Value = '  -2.77%  '
$ws.Range('D21').Value = '13.13'
$ws.Range('E21').Value = '  -4.42%  '
$ws.Range('E22').Value = '  -2.47%  '
$ws.Range('D23').Value = '68.86'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').Value = '262.63'
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('D26').Value = '0.172'
$ws.Range('E26').Value = '  -3.87%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = '25.99'
$ws.Range('E28').Value = '  -4.23%  '
$ws.Range('D29').Value = '7.21'
$ws.Range('E29').Value = '  -5.35%  '
$ws.Range('D30').Value = '6.85'
$ws.Range('E30').Value = '  +6.27%  '
$ws.Range('E31').Value = '  -3.55%  '
$ws.Range('D32').Value = '9.96'
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('D33').Value = '34.99'
$ws.Range('E33').Value = '  -4.87%  '
$ws.Range('E34').Value = '  -3.43%  '
$ws.Range('D35').Value = '50.25'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('D36').Value = '0.0431'
$ws.Range('E36').Value = '  -3.34%  '
$ws.Range('E38').Value = '  -6.62%  '
$ws.Range('D39').Value = '17.27'
$ws.Range('E39').Value = '  -4.36%  '
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('E41').Value = '  -7.55%  '
$ws.Range('E42').Value = '  -2.95%  '
$ws.Range('D43').Value = '22.22'
$ws.Range('E43').Value = '  -3.59%  '
$ws.Range('D44').Value = '119.46'
$ws.Range('E44').Value = '  -3.34%  '
$ws.Range('E45').Value = '  -2.60%  '
$ws.Range('D46').Value = '2.039.00'
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('E48').Value = '  -5.50%  '
$ws.Range('D49').Value = '0.264'
$ws.Range('E49').Value = '  +6.44%  '
$ws.Range('D50').Value = '3.237.34'
$ws.Range('E50').Value = '  -2.30%  '
$ws.Range('E51').Value = '  -3.19%  '
